$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("R40" rule): the rule-name cell B11 is updated from the text
# "R40" to the text "1". Force a Text number format before assigning so
# the numeric-looking literal is stored as a (shared) string rather than
# being auto-coerced into a number.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
